$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - 320 Pitt Street, Sydney, 2000
$ws.Range("D2").Value = "Data\Output\MapScreenshots\320 Pitt Street, Sydney, 2000_210329.035507.jpg"
$ws.Range("E2").Value = "210329.035508_Success"

# Row 3 - 123 Pitt Street, Sydney, NSW
$ws.Range("D3").Value = "Data\Output\MapScreenshots\123 Pitt Street, Sydney, NSW_210329.035519.jpg"
$ws.Range("E3").Value = "210329.035519_Success"

# Row 4 - 555 Anzac Parade Kingsford 2032
$ws.Range("D4").Value = "Data\Output\MapScreenshots\555 Anzac Parade Kingsford 2032_210329.035531.jpg"
$ws.Range("E4").Value = "210329.035531_Success"
